$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.300.28"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "1.551.73"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.479"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "1.774.37"
$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("D13").Value = "1.550.35"
$ws.Range("E13").Value = "  -1.48%  "

$ws.Range("D14").Value = "28.325.72"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("E20").Value = "  -2.44%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.17%  "

$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("E31").Value = "  -4.59%  "

$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").Value = "1.389.17"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("E34").Value = "  -2.90%  "

$ws.Range("E35").Value = "  +2.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.05%  "

$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("E39").Value = "  -3.10%  "

$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("E41").Value = "  +1.62%  "

$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.776"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "

$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("D47").Value = "1.686.82"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.902"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.39%  "

$ws.Range("E51").Value = "  +0.61%  "
